# Replace the "Problem Definition" picture (a screenshotted text image) on
# slide 6 with a live, editable text box containing the same content.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# Locate + remove the old picture shape (Google Shape;81;p18).
$old = $s.Shapes.Item(1)
$old.Delete()

# EMU -> point helper (1 pt = 12700 EMU) so the new box lands exactly where
# the picture used to be (off x=570450,y=1237925 ext cx=8003100,cy=2371200).
$left = 570450 / 12700
$top = 1237925 / 12700
$width = 8003100 / 12700
$height = 2371200 / 12700

$shp = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$shp.Name = "Google Shape;81;p18"

# Match the picture's old "no fill / no outline" look.
$shp.Fill.Visible = $false
$shp.Line.Visible = $false

# bodyPr: square wrap, top anchor, no autofit, 91425 EMU insets on all sides.
$tf = $shp.TextFrame
$tf.WordWrap = -1
$tf.AutoSize = 0
$tf.VerticalAnchor = 1
$ins = 91425 / 12700
$tf.MarginLeft = $ins
$tf.MarginRight = $ins
$tf.MarginTop = $ins
$tf.MarginBottom = $ins

$tr = $tf.TextRange
$tr.Text = "Problem Definition`r`rThis project is all about Fire Extinguishing Robot which has AI built into it.It has an ability to trigger itself accordingly when it senses fire around it.`r`rThe robot is completely autonomous which works on the sensors and using the fire servers which is connected to the fire sensors.`r`rIt will be designed to efficiently extinguish the fire in no time . `r `r"

# Paragraph-by-paragraph formatting, mirroring the source markup:
#  1 "Problem Definition"                              -> centered, 24pt
#  2 ""                                                 -> left, blank spacer
#  3 "This project is all about..."                     -> left
#  4 ""                                                 -> left, blank spacer
#  5 "The robot is completely autonomous..."            -> left
#  6 ""                                                 -> left, blank spacer
#  7 "It will be designed to efficiently..."            -> left
#  8 " "                                                -> left
#  9 ""                                                 -> left, blank spacer
for ($i = 1; $i -le 9; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $pf = $para.ParagraphFormat
    if ($i -eq 1) {
        $pf.Alignment = 2
        $para.Font.Size = 24
    } else {
        $pf.Alignment = 1
    }
    $pf.SpaceBefore = 0
    $pf.SpaceAfter = 0
    $pf.Bullet.Visible = $false
}

Write-Output ("Shapes on slide 6: " + $s.Shapes.Count)
Write-Output ("New shape text: " + $tr.Text)
